# "update project to UE5"
# GameData.xlsx / Actors sheet: the Cube actor's LocationZ (column B, row 2)
# moves from 200 to 400, and the sheet's saved cursor position moves from
# I9 to B3 (reflecting where the author was last working in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actors")

# LocationZ for "Cube" (row 2): 200 -> 400
$ws.Range("B2").Value = 400

# Leave the selection on B3, matching the saved sheetView selection.
$ws.Range("B3").Select()
